$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.02677
$ws.Cells.Item(2, 4).Value = 1.144
$ws.Cells.Item(2, 5).Value = 2.38
$ws.Cells.Item(2, 6).Value = 5
$ws.Cells.Item(3, 3).Value = 0.07140000000000001
$ws.Cells.Item(3, 4).Value = 1.416
$ws.Cells.Item(3, 5).Value = 1.819
$ws.Cells.Item(3, 6).Value = 5
$ws.Cells.Item(4, 3).Value = 0.0346
$ws.Cells.Item(4, 4).Value = 1.884
$ws.Cells.Item(4, 5).Value = 1.742
$ws.Cells.Item(4, 6).Value = 5
$ws.Cells.Item(5, 3).Value = 0.0005845
$ws.Cells.Item(5, 4).Value = 2.189
$ws.Cells.Item(5, 5).Value = 1.121
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(6, 3).Value = 0.008944000000000001
$ws.Cells.Item(6, 4).Value = 2.462
$ws.Cells.Item(6, 5).Value = 0.4541
$ws.Cells.Item(6, 6).Value = 5
$ws.Cells.Item(7, 3).Value = 0.002097
$ws.Cells.Item(7, 4).Value = 2.708
$ws.Cells.Item(7, 5).Value = 0.9320000000000001
$ws.Cells.Item(7, 6).Value = 5
$ws.Cells.Item(8, 3).Value = 0.00004872
$ws.Cells.Item(8, 4).Value = 2.032
$ws.Cells.Item(8, 5).Value = 1.225
$ws.Cells.Item(8, 6).Value = 5
$ws.Cells.Item(9, 3).Value = 0.006281
$ws.Cells.Item(9, 4).Value = 2.769
$ws.Cells.Item(9, 5).Value = 1.321
$ws.Cells.Item(9, 6).Value = 5
$ws.Cells.Item(10, 3).Value = 0.04933
$ws.Cells.Item(10, 4).Value = 3.006
$ws.Cells.Item(10, 5).Value = 0.9503
$ws.Cells.Item(10, 6).Value = 5
$ws.Cells.Item(11, 3).Value = 0.04533
$ws.Cells.Item(11, 4).Value = 2.871
$ws.Cells.Item(11, 5).Value = 0.8944
$ws.Cells.Item(11, 6).Value = 5
$ws.Cells.Item(12, 3).Value = 0.001937
$ws.Cells.Item(12, 4).Value = 2.91
$ws.Cells.Item(12, 5).Value = 0.6382
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(13, 3).Value = 0.001944
$ws.Cells.Item(13, 4).Value = 2.898
$ws.Cells.Item(13, 5).Value = 0.5722
$ws.Cells.Item(13, 6).Value = 5
$ws.Cells.Item(14, 3).Value = 0.2761
$ws.Cells.Item(14, 4).Value = 2.331
$ws.Cells.Item(14, 5).Value = 0.2538
$ws.Cells.Item(14, 6).Value = 15
$ws.Cells.Item(15, 3).Value = 0.1313
$ws.Cells.Item(15, 4).Value = 2.487
$ws.Cells.Item(15, 5).Value = 0.2041
$ws.Cells.Item(15, 6).Value = 10
$ws.Cells.Item(16, 3).Value = 0.0005012
$ws.Cells.Item(16, 4).Value = 3.475
$ws.Cells.Item(16, 5).Value = 0.4351
$ws.Cells.Item(16, 6).Value = 5
$ws.Cells.Item(17, 3).Value = 0.0000004638
$ws.Cells.Item(17, 4).Value = 3.888
$ws.Cells.Item(17, 5).Value = 0.844
$ws.Cells.Item(17, 6).Value = 5
$ws.Cells.Item(18, 3).Value = 0.00008296
$ws.Cells.Item(18, 4).Value = 4.028
$ws.Cells.Item(18, 5).Value = 0.9012
$ws.Cells.Item(18, 6).Value = 5
$ws.Cells.Item(19, 3).Value = 0.02765
$ws.Cells.Item(19, 4).Value = 2.943
$ws.Cells.Item(19, 5).Value = 0.5685
$ws.Cells.Item(19, 6).Value = 10
$ws.Cells.Item(20, 3).Value = 0.0007654
$ws.Cells.Item(20, 4).Value = 4.617
$ws.Cells.Item(20, 5).Value = 1.179
$ws.Cells.Item(20, 6).Value = 5
$ws.Cells.Item(21, 3).Value = 0.0006609
$ws.Cells.Item(21, 4).Value = 4.968
$ws.Cells.Item(21, 5).Value = 1.323
$ws.Cells.Item(21, 6).Value = 5
$ws.Cells.Item(22, 3).Value = 0.007643
$ws.Cells.Item(22, 4).Value = 5.16
$ws.Cells.Item(22, 5).Value = 1.304
$ws.Cells.Item(22, 6).Value = 5
$ws.Cells.Item(23, 3).Value = 0.003798
$ws.Cells.Item(23, 4).Value = 5.25
$ws.Cells.Item(23, 5).Value = 0.7727000000000001
$ws.Cells.Item(23, 6).Value = 5
$ws.Cells.Item(24, 3).Value = 0.0002045
$ws.Cells.Item(24, 4).Value = 5.529
$ws.Cells.Item(24, 5).Value = 0.5427
$ws.Cells.Item(24, 6).Value = 5
$ws.Cells.Item(25, 3).Value = 0.03006
$ws.Cells.Item(25, 4).Value = 4.23
$ws.Cells.Item(25, 5).Value = 0.05253
$ws.Cells.Item(25, 6).Value = 10
$ws.Cells.Item(26, 3).Value = 0.03846
$ws.Cells.Item(26, 4).Value = 4.002
$ws.Cells.Item(26, 5).Value = 0.146
$ws.Cells.Item(26, 6).Value = 15
$ws.Cells.Item(27, 3).Value = 0.0008876999999999999
$ws.Cells.Item(27, 4).Value = 5.346
$ws.Cells.Item(27, 5).Value = 0.4769
$ws.Cells.Item(27, 6).Value = 5
$ws.Cells.Item(28, 3).Value = 0.00139
$ws.Cells.Item(28, 4).Value = 5.644
$ws.Cells.Item(28, 5).Value = 0.8708
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(29, 3).Value = 0.001937
$ws.Cells.Item(29, 4).Value = 5.772
$ws.Cells.Item(29, 5).Value = 0.9781
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(30, 3).Value = 0.0002671
$ws.Cells.Item(30, 4).Value = 5.768
$ws.Cells.Item(30, 5).Value = 1.092
$ws.Cells.Item(30, 6).Value = 5
$ws.Cells.Item(31, 3).Value = 0.0001061
$ws.Cells.Item(31, 4).Value = 5.86
$ws.Cells.Item(31, 5).Value = 0.8046
$ws.Cells.Item(31, 6).Value = 5
$ws.Cells.Item(32, 3).Value = 0.0002671
$ws.Cells.Item(32, 4).Value = 5.846
$ws.Cells.Item(32, 5).Value = 0.2972
$ws.Cells.Item(32, 6).Value = 5
$ws.Cells.Item(33, 3).Value = 0.0004712
$ws.Cells.Item(33, 4).Value = 5.749
$ws.Cells.Item(33, 5).Value = 0.4281
$ws.Cells.Item(33, 6).Value = 5
$ws.Cells.Item(34, 3).Value = 0.02777
$ws.Cells.Item(34, 4).Value = 5.092
$ws.Cells.Item(34, 5).Value = 0.0883
$ws.Cells.Item(34, 6).Value = 15
$ws.Cells.Item(35, 3).Value = 0.001494
$ws.Cells.Item(35, 4).Value = 5.901
$ws.Cells.Item(35, 5).Value = 0.4095
$ws.Cells.Item(35, 6).Value = 5
$ws.Cells.Item(36, 3).Value = 0.002922
$ws.Cells.Item(36, 4).Value = 6.127
$ws.Cells.Item(36, 5).Value = 0.8233
$ws.Cells.Item(36, 6).Value = 5
$ws.Cells.Item(37, 3).Value = 0.0006306
$ws.Cells.Item(37, 4).Value = 6.089
$ws.Cells.Item(37, 5).Value = 0.6877
$ws.Cells.Item(37, 6).Value = 5
$ws.Cells.Item(38, 3).Value = 0.0001719
$ws.Cells.Item(38, 4).Value = 6.185
$ws.Cells.Item(38, 5).Value = 0.476
$ws.Cells.Item(38, 6).Value = 5
$ws.Cells.Item(39, 3).Value = 0.001069
$ws.Cells.Item(39, 4).Value = 6.335
$ws.Cells.Item(39, 5).Value = 0.4138
$ws.Cells.Item(39, 6).Value = 5
$ws.Cells.Item(40, 3).Value = 0.0002812
$ws.Cells.Item(40, 4).Value = 6.341
$ws.Cells.Item(40, 5).Value = 0.3649
$ws.Cells.Item(40, 6).Value = 5
$ws.Cells.Item(41, 3).Value = 0.02613
$ws.Cells.Item(41, 4).Value = 5.043
$ws.Cells.Item(41, 5).Value = 0.037
$ws.Cells.Item(41, 6).Value = 15
$ws.Cells.Item(42, 3).Value = 0.000007891
$ws.Cells.Item(42, 4).Value = 5.954
$ws.Cells.Item(42, 5).Value = 0.6072
$ws.Cells.Item(42, 6).Value = 5
$ws.Cells.Item(43, 3).Value = 0.4577
$ws.Cells.Item(43, 4).Value = 5.047
$ws.Cells.Item(43, 5).Value = 0.4592
$ws.Cells.Item(43, 6).Value = 25
$ws.Cells.Item(44, 3).Value = 0.0002045
$ws.Cells.Item(44, 4).Value = 6.582
$ws.Cells.Item(44, 5).Value = 0.5487
$ws.Cells.Item(44, 6).Value = 5
$ws.Cells.Item(45, 3).Value = 0.002354
$ws.Cells.Item(45, 4).Value = 6.634
$ws.Cells.Item(45, 5).Value = 0.9028
$ws.Cells.Item(45, 6).Value = 5
$ws.Cells.Item(46, 3).Value = 0.00129
$ws.Cells.Item(46, 4).Value = 6.625
$ws.Cells.Item(46, 5).Value = 0.8683999999999999
$ws.Cells.Item(46, 6).Value = 5
$ws.Cells.Item(47, 3).Value = 0.00479
$ws.Cells.Item(47, 4).Value = 6.751
$ws.Cells.Item(47, 5).Value = 0.483
$ws.Cells.Item(47, 6).Value = 5
$ws.Cells.Item(48, 3).Value = 0.00003551
$ws.Cells.Item(48, 4).Value = 7.113
$ws.Cells.Item(48, 5).Value = 0.5113
$ws.Cells.Item(48, 6).Value = 5
$ws.Cells.Item(49, 3).Value = 0.00008609
$ws.Cells.Item(49, 4).Value = 7.171
$ws.Cells.Item(49, 5).Value = 0.6042
$ws.Cells.Item(49, 6).Value = 5
$ws.Cells.Item(50, 3).Value = 0.1182
$ws.Cells.Item(50, 4).Value = 5.289
$ws.Cells.Item(50, 5).Value = 0.3804
$ws.Cells.Item(50, 6).Value = 10
$ws.Cells.Item(51, 3).Value = 0.00001403
$ws.Cells.Item(51, 4).Value = 6.289
$ws.Cells.Item(51, 5).Value = 1.102
$ws.Cells.Item(51, 6).Value = 5
$ws.Cells.Item(52, 3).Value = 0.003351
$ws.Cells.Item(52, 4).Value = 6.51
$ws.Cells.Item(52, 5).Value = 1.376
$ws.Cells.Item(52, 6).Value = 5
$ws.Cells.Item(53, 3).Value = 0.0003045
$ws.Cells.Item(53, 4).Value = 6.332
$ws.Cells.Item(53, 5).Value = 1.739
$ws.Cells.Item(53, 6).Value = 5
$ws.Cells.Item(54, 3).Value = 0.001467
$ws.Cells.Item(54, 4).Value = 6.269
$ws.Cells.Item(54, 5).Value = 1.21
$ws.Cells.Item(54, 6).Value = 5
$ws.Cells.Item(55, 3).Value = 0.001435
$ws.Cells.Item(55, 4).Value = 5.811
$ws.Cells.Item(55, 5).Value = 0.9405
$ws.Cells.Item(55, 6).Value = 5
$ws.Cells.Item(56, 3).Value = 0.003602
$ws.Cells.Item(56, 4).Value = 5.425
$ws.Cells.Item(56, 5).Value = 1.215
$ws.Cells.Item(56, 6).Value = 5
$ws.Cells.Item(57, 3).Value = 0.002066
$ws.Cells.Item(57, 4).Value = 5.557
$ws.Cells.Item(57, 5).Value = 0.4551
$ws.Cells.Item(57, 6).Value = 5
$ws.Cells.Item(58, 3).Value = 0.001003
$ws.Cells.Item(58, 4).Value = 5.238
$ws.Cells.Item(58, 5).Value = 0.4266
$ws.Cells.Item(58, 6).Value = 5
$ws.Cells.Item(59, 3).Value = 0.008064999999999999
$ws.Cells.Item(59, 4).Value = 5.144
$ws.Cells.Item(59, 5).Value = 0.5847
$ws.Cells.Item(59, 6).Value = 5
$ws.Cells.Item(60, 3).Value = 0.006485
$ws.Cells.Item(60, 4).Value = 5.281
$ws.Cells.Item(60, 5).Value = 0.1414
$ws.Cells.Item(60, 6).Value = 5
$ws.Cells.Item(61, 3).Value = 0.2548
$ws.Cells.Item(61, 4).Value = 5.083
$ws.Cells.Item(61, 5).Value = 0.191
$ws.Cells.Item(61, 6).Value = 15
$ws.Cells.Item(62, 3).Value = 0.02748
$ws.Cells.Item(62, 4).Value = 4.555
$ws.Cells.Item(62, 5).Value = 0.5453
$ws.Cells.Item(62, 6).Value = 10
$ws.Cells.Item(63, 3).Value = 0.001621
$ws.Cells.Item(63, 4).Value = 5.397
$ws.Cells.Item(63, 5).Value = 1.361
$ws.Cells.Item(63, 6).Value = 5
$ws.Cells.Item(64, 3).Value = 0.0003016
$ws.Cells.Item(64, 4).Value = 5.459
$ws.Cells.Item(64, 5).Value = 1.516
$ws.Cells.Item(64, 6).Value = 5
$ws.Cells.Item(65, 3).Value = 0.001136
$ws.Cells.Item(65, 4).Value = 5.798
$ws.Cells.Item(65, 5).Value = 1.697
$ws.Cells.Item(65, 6).Value = 5
$ws.Cells.Item(66, 3).Value = 0.004142
$ws.Cells.Item(66, 4).Value = 6.247
$ws.Cells.Item(66, 5).Value = 2.258
$ws.Cells.Item(66, 6).Value = 5
$ws.Cells.Item(67, 3).Value = 0.0009496
$ws.Cells.Item(67, 4).Value = 6.376
$ws.Cells.Item(67, 5).Value = 2.06
$ws.Cells.Item(67, 6).Value = 5
$ws.Cells.Item(68, 3).Value = 0.001266
$ws.Cells.Item(68, 4).Value = 6.621
$ws.Cells.Item(68, 5).Value = 1.491
$ws.Cells.Item(68, 6).Value = 5
$ws.Cells.Item(69, 3).Value = 0.002644
$ws.Cells.Item(69, 4).Value = 6.99
$ws.Cells.Item(69, 5).Value = 1.061
$ws.Cells.Item(69, 6).Value = 5
$ws.Cells.Item(70, 3).Value = 0.000004899
$ws.Cells.Item(70, 4).Value = 7.357
$ws.Cells.Item(70, 5).Value = 0.9724
$ws.Cells.Item(70, 6).Value = 5
$ws.Cells.Item(71, 3).Value = 0.0002169
$ws.Cells.Item(71, 4).Value = 7.45
$ws.Cells.Item(71, 5).Value = 0.3644
$ws.Cells.Item(71, 6).Value = 5

# Add new row 72 (extend the table by one record)
$ws.Range("A71").Copy()
$ws.Range("A72").PasteSpecial(-4122)
$ws.Cells.Item(72, 1).Value = 70
$ws.Cells.Item(72, 2).Value = 71
$ws.Cells.Item(72, 3).Value = 0.6113
$ws.Cells.Item(72, 4).Value = 5.114
$ws.Cells.Item(72, 5).Value = 0.1405
$ws.Cells.Item(72, 6).Value = 25
